$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: A2/B2 numeric changes
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 7161

# Row 2: text changes (was "ELEGANCE LINE SHADOW POWDER 4")
$ws.Range("C2").Value = "ELEGANCE LIP CONTOUR 2044 ORANGE"
$ws.Range("D2").Value = "ELEGANCE LIP CONTOUR 2044 ORANGE"
$ws.Range("I2").Value = "ELEGANCE LIP CONTOUR 2044 ORANGE"
$ws.Range("J2").Value = "ELEGANCE LIP CONTOUR 2044 ORANGE"

# Row 3: A3/B3 numeric changes
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 7162

# Row 3: text changes (was "ELEGANCE LIP CONTOUR 2035 NUDE BROWN")
$ws.Range("C3").Value = "ELEGANCE LIP CONTOUR BROWN"
$ws.Range("D3").Value = "ELEGANCE LIP CONTOUR BROWN"
$ws.Range("I3").Value = "ELEGANCE LIP CONTOUR BROWN"
$ws.Range("J3").Value = "ELEGANCE LIP CONTOUR BROWN"
